$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the dataset. It belongs right before the
# existing row 88, so insert a new row there and push everything else down
# (row 88 -> 89, ..., row 100 -> 101).
$ws.Rows("88:88").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A88").Value = 6
$ws.Range("B88").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C88").Value = "Metropolitana"
$ws.Range("D88").Value = 45132
$ws.Range("E88").Value = 13
$ws.Range("F88").Value = 100112035
$ws.Range("G88").Value = "Bruselas (repollito)"
$ws.Range("H88").Value = "Sin especificar"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 290
$ws.Range("K88").Value = 18000
$ws.Range("L88").Value = 20000
$ws.Range("M88").Value = 19172
$ws.Range("N88").Value = '$/malla 15 kilos'
$ws.Range("O88").Value = "Provincia de Quillota"
$ws.Range("P88").Value = 1278
$ws.Range("Q88").Value = 15
$ws.Range("R88").Value = "Hortaliza"

# Make sure the date column keeps the same date/time number format used by
# the rest of column D.
$ws.Range("D88").NumberFormat = $ws.Range("D89").NumberFormat
